$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Database table names now require the explicit "project." schema
#    prefix. Update the "Tables" sheet's table-name column (B) for
#    every affected table definition.
# ------------------------------------------------------------------
$tables = $wb.Worksheets.Item("Tables")

$tables.Range("B2").Value = "project.site"    # hidden.available_sites
$tables.Range("B3").Value = "project.device"  # hidden.available_systems
$tables.Range("B4").Value = "project.site"    # hidden.lease_boundaries
$tables.Range("B5").Value = "project.site"    # hidden.site_boundaries
$tables.Range("B6").Value = "project.site"    # hidden.corridor_boundaries
$tables.Range("B7").Value = "project.farm"    # hidden.landing_points

# ------------------------------------------------------------------
# 2) Make "Tables" the active/selected sheet (was "ROOT"), with B7
#    selected as the active cell.
# ------------------------------------------------------------------
$tables.Activate()
$tables.Range("B7").Select()

# ------------------------------------------------------------------
# 3) The ROOT sheet's header row no longer carries an explicit custom
#    row height - let it revert to the sheet's default height.
# ------------------------------------------------------------------
$root = $wb.Worksheets.Item("ROOT")
$root.Rows.Item(1).AutoFit()
